$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.282.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.00%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.409.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.35%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.78%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.408.88'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.33%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.495'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.53%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.19'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.35%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.120'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -10.27%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.373'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.79%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.993.02'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000178'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -11.22%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.421.08'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.88%  '

$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.115'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.52%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.301.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.94%  '

$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.95'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -9.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -11.28%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.84'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.54'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.82%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '385.64'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.553'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.65%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.06%  '

$ws.Range("E25").Value = '  +0.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.553.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000105'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -10.15%  '

$ws.Range("E28").Value = '  +0.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.06'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.87%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.08'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.96%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.20'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -9.86%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.419.31'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.00%  '

$ws.Range("E33").Value = '  -0.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.143'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '22.78'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.58%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '171.29'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.75'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -10.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.13'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -12.53%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.45'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.68%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.71'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -10.76%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0757'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.71%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.813'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.95%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -13.68%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.59'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -11.32%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.39%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.15'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.30%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.49'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.45%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.04'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -15.97%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.177.96'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.22%  '
